$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch cell (far outside the used range) formatted as Text to hold each
# new price value, then paste-special (values only) into the target D-column cell.
# This keeps the destination cell a plain text cell (matching the original
# inlineStr/shared-string "Price" column) instead of Excel auto-converting the
# numeric-looking text into a real number.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$scratch.Value = "266.29"
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$scratch.Value = "21.38"
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)

$scratch.Value = "6.172"
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)

$scratch.Value = "0.06161"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)

$scratch.Value = "3.568"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)

$scratch.Value = "6.521"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)

$scratch.Value = "0.8249"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)

$scratch.Value = "0.01348"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)

$scratch.Value = "0.1584"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)

$scratch.Value = "0.08121"
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)

$scratch.Value = "0.03348"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)

$scratch.Value = "0.03182"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)

$scratch.Value = "0.09255"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)

$scratch.Value = "3.769"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)

$scratch.Value = "0.001629"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)

$scratch.Value = "0.04677"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)

$scratch.Value = "0.006403"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)

$scratch.Value = "0.006203"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)

$scratch.Value = "0.0001498"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)

$scratch.Value = "3.734"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)

$scratch.Value = "2.431"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)

$scratch.Value = "0.3300"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)

$scratch.Value = "0.1239"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)

$scratch.Value = "0.04652"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)

$scratch.Value = "0.006973"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)

$scratch.Value = "0.1125"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)

$scratch.Value = "0.003658"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)

$scratch.Value = "0.01158"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)

$scratch.Value = "0.00005941"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)

$scratch.Value = "0.0009877"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)

$scratch.Value = "0.00000000749"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)

$scratch.Value = "0.002439"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)

$scratch.Value = "0.00001898"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)

$scratch.Value = "0.01238"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)

# Clean up the scratch cell so it leaves no trace in the saved workbook.
$scratch.Clear()
$excel.CutCopyMode = 0
